$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.01253208636536152
$ws.Range("C2").Value = 0.04103571897497393
$ws.Range("D2").Value = 186123.597850132
$ws.Range("E2").Value = 2797.565817734744
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 188921.2172356721
